$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D3").Value = -7.857000000000001
$ws.Range("D4").Value = -8.074999999999999
$ws.Range("D7").Value = -7.935
$ws.Range("D8").Value = -8.108000000000001
$ws.Range("C11").Value = -12.843
$ws.Range("C12").Value = -13.117
$ws.Range("D12").Value = -7.992999999999999
$ws.Range("D14").Value = -8.178999999999998
$ws.Range("C15").Value = -12.18
$ws.Range("D22").Value = -7.812
